$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted above the existing row 111,
# pushing every subsequent row (old 111..161) down by one (new 112..162).
$ws.Rows(111).Insert()

# Populate the newly inserted row 111 with its own data.
$ws.Range("A111").Value = 3
$ws.Range("B111").Value = "Femacal de La Calera"
$ws.Range("C111").Value = "Coquimbo"
$ws.Range("D111").Value = 44489
$ws.Range("E111").Value = 5
$ws.Range("F111").Value = 100112001
$ws.Range("G111").Value = "Berenjena"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 95
$ws.Range("K111").Value = 8000
$ws.Range("L111").Value = 8500
$ws.Range("M111").Value = 8237
$ws.Range("N111").Value = "$/caja 60 unidades"
$ws.Range("O111").Value = "Región de Arica y Parinacota"
$ws.Range("P111").Value = 137
$ws.Range("Q111").Value = 60
$ws.Range("R111").Value = "Hortaliza"
